{"js": "// Update the WMO Met Data Exchange Interoperability Experiment document:\n//  1. Bump the document date (title-page date + the \"Date:\" row in the\n//     metadata table) from 2025-02-13 to 2025-02-18.\n//  2. Add a new \"Document status: DRAFT\" row to the metadata table,\n//     immediately after the \"Version: 0.1\" row.\n\nconst body = context.document.body;\n\n// --- 1. Replace every occurrence of the old date with the new one -------\nconst dateResults = body.search(\"2025-02-13\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < dateResults.items.length; i++) {\n  dateResults.items[i].insertText(\"2025-02-18\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 2. Insert a new \"Document status: DRAFT\" row after \"Version: 0.1\" --\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nrows.items.forEach((r) => r.load(\"values\"));\nawait context.sync();\n\nlet versionRow = null;\nfor (let i = 0; i < rows.items.length; i++) {\n  const values = rows.items[i].values;\n  if (values && values[0] && values[0][0] && values[0][0].indexOf(\"Version:\") === 0) {\n    versionRow = rows.items[i];\n    break;\n  }\n}\n\nif (versionRow) {\n  versionRow.insertRows(Word.InsertLocation.after, 1, [[\"Document status: DRAFT\"]]);\n} else {\n  // Fallback: append at the end of the table if the \"Version:\" row can't be found.\n  table.addRows(Word.InsertLocation.end, 1, [[\"Document status: DRAFT\"]]);\n}\nawait context.sync();\n", "ps1": "# Update the WMO Met Data Exchange Interoperability Experiment document:\n#  1. Bump the document date (title-page date + the \"Date:\" row in the\n#     metadata table) from 2025-02-13 to 2025-02-18.\n#  2. Add a new \"Document status: DRAFT\" row to the metadata table,\n#     immediately after the \"Version: 0.1\" row.\n\n$d = $word.ActiveDocument\n\n# --- 1. Replace every occurrence of the old date with the new one -------\n$find = $d.Content.Find\n$find.Text = \"2025-02-13\"\n$find.Replacement.Text = \"2025-02-18\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# --- 2. Insert a new \"Document status: DRAFT\" row after \"Version: 0.1\" --\n$t = $d.Tables.Item(1)\n\n$versionIndex = 0\nfor ($i = 1; $i -le $t.Rows.Count; $i++) {\n    $cellText = $t.Rows.Item($i).Cells.Item(1).Range.Text\n    if ($cellText -like \"Version:*\") {\n        $versionIndex = $i\n        break\n    }\n}\n\nif ($versionIndex -gt 0) {\n    if ($versionIndex -lt $t.Rows.Count) {\n        $newRow = $t.Rows.Add($t.Rows.Item($versionIndex + 1))\n    } else {\n        $newRow = $t.Rows.Add()\n    }\n    $newRow.Cells.Item(1).Range.Text = \"Document status: DRAFT\"\n}\n"}
